$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: add a new column O (mirroring column N's
#     formatting) with a 0 value for the data rows, and an empty but
#     formatted cell in the opening balance row (row 2). Cell writes/copy
#     operate on the sheet object directly so the active tab isn't
#     disturbed. ---
$repay = $wb.Worksheets.Item("Repayment schedule")

# Row 2 keeps O2 blank, just formatted like its neighbours (M2/N2/P2).
$repay.Range("N2").Copy()
[void]$repay.Range("O2").PasteSpecial(-4122)

# Rows 3-8 get an explicit 0 value, formatted like column N.
foreach ($row in 3..8) {
    $srcCell = "N" + $row
    $dstCell = "O" + $row
    $repay.Range($srcCell).Copy()
    [void]$repay.Range($dstCell).PasteSpecial(-4122)
    $repay.Range($dstCell).Value = 0
}

$excel.CutCopyMode = $false

# --- Summary sheet: move the selection from A7:XFD15 to C3 ---
$summary = $wb.Worksheets.Item("Summary")
[void]$summary.Range("C3").Select()

# --- Edit Repayment Schedule sheet: move the selection from A8 to A9 ---
$editSheet = $wb.Worksheets.Item("Edit Repayment Schedule")
[void]$editSheet.Range("A9").Select()

# --- Restore the originally active sheet/tab (NewLoanInput) ---
[void]$wb.Worksheets.Item("NewLoanInput").Activate()
